# "Generate Report for Handback"
#
# The handback transform for file 5d8d9ab8-3244-4662-82ba-50c2d91030ad failed
# because the handback file name didn't match the handoff file name. Update
# the localization status report to reflect the failure:
#   - Status for that file becomes "Handback transform failed" (was
#     "Ready for handoff") on the Overview sheet and on each language sheet.
#   - The "Error Detail" column (P) for that file's row on the zh-cn and
#     de-de sheets gets the explanatory message.
#   - The Error Detail column is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 7 on every sheet is the 5d8d9ab8-3244-4662-82ba-50c2d91030ad.md file.
# Its status flips from "Ready for handoff" to "Handback transform failed".
$wsOverview.Range("E7").Value = "Handback transform failed"
$wsOverview.Range("F7").Value = "Handback transform failed"
$wsZhCn.Range("C7").Value = "Handback transform failed"
$wsDeDe.Range("C7").Value = "Handback transform failed"

# Record why the handback transform failed, per language, in the
# "Error Detail" column (P) of that same row.
$wsZhCn.Range("P7").Value = "Handback file name: atiplge4.doj is different with handoff file name: 5d8d9ab8-3244-4662-82ba-50c2d91030ad.c4dc661659e8a39bb160aec482287a6ee4d86b8d.zh-cn."
$wsDeDe.Range("P7").Value = "Handback file name: atiplge4.doj is different with handoff file name: 5d8d9ab8-3244-4662-82ba-50c2d91030ad.c4dc661659e8a39bb160aec482287a6ee4d86b8d.de-de."

# Widen the Error Detail column so the new message is readable
# (ColumnWidth 39.17 Excel units serializes to the OOXML width of 40).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
